# Fruta / hortaliza, semanal
# Re-shuffle the per-row data (Fecha, Calidad, Volumen, Precio min/max/prom,
# Origen, Precio $/Kg) across rows 2-37 according to the new weekly snapshot.
# A, B, C, E, F, G, H, N, Q, R are constant across the whole table and are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "I", "J", "K", "L", "M", "O", "P")

# Row(destination) -> Row(source of the snapshot values), a permutation of
# rows 2..37.
$map = @{
    2  = 33
    3  = 23
    4  = 27
    5  = 24
    6  = 10
    7  = 28
    8  = 14
    9  = 2
    10 = 29
    11 = 35
    12 = 34
    13 = 13
    14 = 11
    15 = 37
    16 = 20
    17 = 16
    18 = 17
    19 = 19
    20 = 18
    21 = 4
    22 = 5
    23 = 36
    24 = 32
    25 = 12
    26 = 6
    27 = 31
    28 = 26
    29 = 15
    30 = 8
    31 = 21
    32 = 3
    33 = 22
    34 = 9
    35 = 7
    36 = 25
    37 = 30
}

# 1) Snapshot every source cell's current value before any write happens -
#    the mapping is a permutation, not a simple shift, so rows get both
#    read from and written to.
$snapshot = @{}
foreach ($col in $cols) {
    $snapshot[$col] = @{}
    for ($row = 2; $row -le 37; $row++) {
        $snapshot[$col][$row] = $ws.Range($col + $row).Value2
    }
}

# 2) Write the snapshot back out according to the row mapping.
foreach ($col in $cols) {
    foreach ($destRow in $map.Keys) {
        $srcRow = $map[$destRow]
        $ws.Range($col + $destRow).Value = $snapshot[$col][$srcRow]
    }
}
